# Updated cryptos list on Thu Oct 19 08:57:03 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row = 2; D = '28.341.09'; E = '  -0.85%  ' }
    @{ Row = 3; D = '1.550.21'; E = '  -1.79%  ' }
    @{ Row = 4; D = $null; E = '  -0.11%  ' }
    @{ Row = 5; D = '209.90'; E = '  -1.59%  ' }
    @{ Row = 6; D = '0.481'; E = '  -1.71%  ' }
    @{ Row = 7; D = $null; E = '  -0.17%  ' }
    @{ Row = 8; D = '23.84'; E = '  -0.96%  ' }
    @{ Row = 9; D = $null; E = '  -1.84%  ' }
    @{ Row = 11; D = '0.0888'; E = '  -0.35%  ' }
    @{ Row = 12; D = '1.772.55'; E = '  -1.72%  ' }
    @{ Row = 13; D = '1.543.97'; E = '  -1.37%  ' }
    @{ Row = 14; D = '28.317.12'; E = $null }
    @{ Row = 15; D = $null; E = '  -2.05%  ' }
    @{ Row = 16; D = $null; E = '  -1.89%  ' }
    @{ Row = 17; D = '60.75'; E = '  -2.28%  ' }
    @{ Row = 18; D = '227.48'; E = '  -1.77%  ' }
    @{ Row = 19; D = $null; E = '  -0.63%  ' }
    @{ Row = 20; D = $null; E = '  -2.39%  ' }
    @{ Row = 21; D = $null; E = '  -0.14%  ' }
    @{ Row = 22; D = '3.90'; E = '  +0.19%  ' }
    @{ Row = 23; D = $null; E = '  -2.72%  ' }
    @{ Row = 24; D = '2.02'; E = '  -1.95%  ' }
    @{ Row = 25; D = '151.31'; E = '  +0.08%  ' }
    @{ Row = 26; D = $null; E = '  -2.01%  ' }
    @{ Row = 27; D = $null; E = '  -1.14%  ' }
    @{ Row = 28; D = $null; E = '  -0.15%  ' }
    @{ Row = 29; D = $null; E = '  -3.09%  ' }
    @{ Row = 30; D = '0.0466'; E = '  -3.80%  ' }
    @{ Row = 31; D = $null; E = '  -4.75%  ' }
    @{ Row = 32; D = $null; E = '  -1.44%  ' }
    @{ Row = 33; D = '1.383.40'; E = '  -1.14%  ' }
    @{ Row = 34; D = $null; E = '  -3.25%  ' }
    @{ Row = 35; D = $null; E = '  +1.58%  ' }
    @{ Row = 36; D = $null; E = '  -3.58%  ' }
    @{ Row = 37; D = $null; E = '  -1.13%  ' }
    @{ Row = 38; D = $null; E = '  -2.06%  ' }
    @{ Row = 39; D = $null; E = '  -2.79%  ' }
    @{ Row = 40; D = $null; E = '  +1.38%  ' }
    @{ Row = 41; D = '0.508'; E = '  -2.53%  ' }
    @{ Row = 42; D = '0.999'; E = '  -0.18%  ' }
    @{ Row = 43; D = '0.777'; E = '  -2.18%  ' }
    @{ Row = 44; D = $null; E = '  -2.17%  ' }
    @{ Row = 45; D = '5.38'; E = '  -1.93%  ' }
    @{ Row = 46; D = $null; E = '  -2.05%  ' }
    @{ Row = 47; D = '1.684.39'; E = '  -1.79%  ' }
    @{ Row = 48; D = '0.866'; E = '  -10.04%  ' }
    @{ Row = 49; D = '85.43'; E = '  -1.34%  ' }
    @{ Row = 50; D = '42.25'; E = '  +5.15%  ' }
    @{ Row = 51; D = $null; E = '  -0.20%  ' }
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($null -ne $u.D) {
        $cellD = $ws.Cells.Item($row, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.D
    }
    if ($null -ne $u.E) {
        $cellE = $ws.Cells.Item($row, 5)
        $cellE.NumberFormat = "@"
        $cellE.Value = $u.E
    }
}
